$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Fatalities by Age Group
# ---------------------------------------------------------------------------
$wsAge = $wb.Worksheets.Item("Fatalities by Age Group")
$wsAge.Range("B2").Value = 12
$wsAge.Range("B4").Value = 59
$wsAge.Range("B5").Value = 506
$wsAge.Range("B6").Value = 1639
$wsAge.Range("B7").Value = 4161
$wsAge.Range("B8").Value = 8180
$wsAge.Range("B9").Value = 6369
$wsAge.Range("B10").Value = 7702
$wsAge.Range("B11").Value = 8350
$wsAge.Range("B12").Value = 8005
$wsAge.Range("B13").Value = 19304

# ---------------------------------------------------------------------------
# Fatalities by Gender
# ---------------------------------------------------------------------------
$wsGender = $wb.Worksheets.Item("Fatalities by Gender")
$wsGender.Range("B2").Value = 26859
$wsGender.Range("B3").Value = 37444

# ---------------------------------------------------------------------------
# Fatalities by Race-Ethnicity
# ---------------------------------------------------------------------------
$wsRace = $wb.Worksheets.Item("Fatalities by Race-Ethnicity")
$wsRace.Range("B2").Value = 1235
$wsRace.Range("B3").Value = 6623
$wsRace.Range("B4").Value = 28476
$wsRace.Range("B5").Value = 365
$wsRace.Range("B6").Value = 27569
$wsRace.Range("B7").Value = 36

# ---------------------------------------------------------------------------
# Restore the view/selection state recorded for each sheet, and make
# "Fatalities by Age Group" the active tab (it was "Fatalities by
# Race-Ethnicity" before).
# ---------------------------------------------------------------------------
$wsRace.Activate()
$wsRace.Range("C12:I20").Select()

$wsGender.Activate()
$wsGender.Range("G12").Select()

$wsAge.Activate()
$wsAge.Range("G5").Select()

# ---------------------------------------------------------------------------
# Main application window geometry, as recorded in the workbook view.
# ---------------------------------------------------------------------------
$win = $excel.ActiveWindow
$win.Left = -96
$win.Top = -96
$win.Width = 23232
$win.Height = 12552
